$wb = $excel.ActiveWorkbook

# The diff adds a new "score" column (E) to the "Option" worksheet: a
# header in E1, plus a numeric score for each option row belonging to
# question q1/q2/q4 (rows 2-12). Rows 13-18 (q5/q5_1 options) are left
# untouched, matching the source diff.
$ws = $wb.Worksheets.Item("Option")
$wsQ = $wb.Worksheets.Item("Question")

# --- header cell E1: same text/number style as the other header cells (D1) ---
$ws.Cells.Item(1, 5).Value = "score"
$ws.Range("D1").Copy()
$ws.Cells.Item(1, 5).PasteSpecial(-4122)   # xlPasteFormats

# --- score values for rows 2-12 ---
$scores = [ordered]@{
    2  = 1.0
    3  = 1.0
    4  = 0.0
    5  = 1.0
    6  = 1.0
    7  = 0.0
    8  = 2.0
    9  = 1.0
    10 = 0.0
    11 = 1.0
    12 = 0.0
}

# Those data cells use the plain, unshaded style (same as used for the
# "Question" sheet's body cells, e.g. A2) rather than the shaded option
# style used in columns A-D.
$wsQ.Range("A2").Copy()
foreach ($row in $scores.Keys) {
    $ws.Cells.Item($row, 5).Value = $scores[$row]
    $ws.Cells.Item($row, 5).PasteSpecial(-4122)   # xlPasteFormats
}
